# Insert a new observation row at row 229 (pushing the existing rows
# 229-257 down to 230-258) for the "Poroto verde" price sheet.
#
# The new row duplicates every field of the (now pushed-down) row 230
# except the date (column D) and the volume (column J), which get their
# own new values. This mirrors the author's edit, which shows row 229
# changing to a brand-new date/volume pair while every row below shifts
# down by one and the previous last row (257) is duplicated into a new
# row 258.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 229..257 down to 230..258, leaving row 229 blank (but with
# formatting inherited from the row above, matching Excel's default
# Insert behaviour).
$ws.Rows.Item(229).Insert()

# Copy every column from the row that used to be 229 (now at 230) into
# the new blank row 229, except D (Fecha) and J (Volumen) which take the
# new values below.
for ($col = 1; $col -le 18; $col++) {
    if ($col -ne 4 -and $col -ne 10) {
        $ws.Cells.Item(229, $col).Value = $ws.Cells.Item(230, $col).Value()
    }
}

# New date (Fecha) and volume (Volumen) for the inserted row.
$ws.Cells.Item(229, 4).Value = 44776
$ws.Cells.Item(229, 10).Value = 440
